$p = $ppt.ActivePresentation

# --- 1) Slide master date placeholder field: 7/1/2019 -> 7/2/2019 -----------
$dateShape = $p.SlideMaster.Shapes.Item("Date Placeholder 3")
$dateRange = $dateShape.TextFrame.TextRange
$dateChars = $dateRange.Characters(1, $dateRange.Length)
$dateChars.Text = "7/2/2019"

# --- 2) Slide 1, "Text Placeholder 2": drop the cached normAutofit line-
#        spacing reduction (lnSpcReduction="10000" -> none) --------------
$s = $p.Slides.Item(1)
$bodyShape = $s.Shapes.Item("Text Placeholder 2")
$bodyShape.TextFrame.AutoSize = 2

# --- 3) Same shape: "android event/tweet" -> "android tweet" ---------------
$bodyRange = $bodyShape.TextFrame.TextRange
$fullText = $bodyRange.Text
$idx = $fullText.IndexOf("event/tweet")
if ($idx -ge 0) {
    $evtTweet = $bodyRange.Characters($idx + 1, 11)
    $evtTweet.Text = "tweet"
}

# --- 4) Slide 1 title "Title 4": merge the three runs into a single run ----
$titleShape = $s.Shapes.Item("Title 4")
$titleRange = $titleShape.TextFrame.TextRange
$titleChars = $titleRange.Characters(1, $titleRange.Length)
$titleChars.Text = "IPhone Users reaction towards Android Tweets"
